$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 3 now-unused trailing rows (21 data rows before -> 18 after: three "group header" rows are folded into column J)
$ws.Range("A19:K21").EntireRow.Delete()

# Force text (string) storage for every data cell we are about to write,
# so numeric-looking strings ("6", "80,000", "5.70", "$8.10", ...) are not
# auto-coerced into Excel numbers/currency by the COM value setter.
$textRange = $ws.Range("A2:K18")
$textRange.NumberFormat = "@"

# Row 1: new numeric column-index header (0..10)
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10

# Row 2: former header row 1 (column captions)
$ws.Range("A2").Value = 'Lg., mm'
$ws.Range("B2").Value = 'Threading'
$ws.Range("C2").Value = 'HeadDia., mm'
$ws.Range("D2").Value = 'HeadHt., mm'
$ws.Range("E2").Value = 'DriveSize, mm'
$ws.Range("F2").Value = 'TensileStrength, psi'
$ws.Range("G2").Value = 'Pkg.Qty.'
$ws.Range("I2").Value = 'Pkg.'
# K2 held the old group-header's material text ("18-8 Stainless Steel"); the new
# header row has nothing there, so explicitly blank it out.
$ws.Range("K2").ClearContents()

# Rows 3-18: data rows (shifted up from original; thread_size group label now filled into column J per row)
$ws.Range("A3").Value = '6'
$ws.Range("B3").Value = 'Fully Threaded'
$ws.Range("C3").Value = '5.70'
$ws.Range("D3").Value = '1.65'
$ws.Range("E3").Value = '2'
$ws.Range("F3").Value = '80,000'
$ws.Range("G3").Value = '5'
$ws.Range("H3").Value = '95942A111'
$ws.Range("I3").Value = '$8.10'
$ws.Range("J3").Value = 'M3 × 0.5 mm'
$ws.Range("K3").Value = '18-8 Stainless Steel'
$ws.Range("A4").Value = '8'
$ws.Range("B4").Value = 'Fully Threaded'
$ws.Range("C4").Value = '5.70'
$ws.Range("D4").Value = '1.65'
$ws.Range("E4").Value = '2'
$ws.Range("F4").Value = '80,000'
$ws.Range("G4").Value = '5'
$ws.Range("H4").Value = '95942A112'
$ws.Range("I4").Value = '8.76'
$ws.Range("J4").Value = 'M3 × 0.5 mm'
$ws.Range("K4").Value = '18-8 Stainless Steel'
$ws.Range("A5").Value = '10'
$ws.Range("B5").Value = 'Fully Threaded'
$ws.Range("C5").Value = '5.70'
$ws.Range("D5").Value = '1.65'
$ws.Range("E5").Value = '2'
$ws.Range("F5").Value = '80,000'
$ws.Range("G5").Value = '5'
$ws.Range("H5").Value = '95942A113'
$ws.Range("I5").Value = '9.54'
$ws.Range("J5").Value = 'M3 × 0.5 mm'
$ws.Range("K5").Value = '18-8 Stainless Steel'
$ws.Range("A6").Value = '12'
$ws.Range("B6").Value = 'Fully Threaded'
$ws.Range("C6").Value = '5.70'
$ws.Range("D6").Value = '1.65'
$ws.Range("E6").Value = '2'
$ws.Range("F6").Value = '80,000'
$ws.Range("G6").Value = '5'
$ws.Range("H6").Value = '95942A114'
$ws.Range("I6").Value = '10.33'
$ws.Range("J6").Value = 'M3 × 0.5 mm'
$ws.Range("K6").Value = '18-8 Stainless Steel'
$ws.Range("A7").Value = '8'
$ws.Range("B7").Value = 'Fully Threaded'
$ws.Range("C7").Value = '7.60'
$ws.Range("D7").Value = '2.20'
$ws.Range("E7").Value = '2.5'
$ws.Range("F7").Value = '80,000'
$ws.Range("G7").Value = '5'
$ws.Range("H7").Value = '95942A115'
$ws.Range("I7").Value = '9.28'
$ws.Range("J7").Value = 'M4 × 0.7 mm'
$ws.Range("K7").Value = '18-8 Stainless Steel'
$ws.Range("A8").Value = '10'
$ws.Range("B8").Value = 'Fully Threaded'
$ws.Range("C8").Value = '7.60'
$ws.Range("D8").Value = '2.20'
$ws.Range("E8").Value = '2.5'
$ws.Range("F8").Value = '80,000'
$ws.Range("G8").Value = '5'
$ws.Range("H8").Value = '95942A116'
$ws.Range("I8").Value = '10.33'
$ws.Range("J8").Value = 'M4 × 0.7 mm'
$ws.Range("K8").Value = '18-8 Stainless Steel'
$ws.Range("A9").Value = '12'
$ws.Range("B9").Value = 'Fully Threaded'
$ws.Range("C9").Value = '7.60'
$ws.Range("D9").Value = '2.20'
$ws.Range("E9").Value = '2.5'
$ws.Range("F9").Value = '80,000'
$ws.Range("G9").Value = '5'
$ws.Range("H9").Value = '95942A117'
$ws.Range("I9").Value = '10.98'
$ws.Range("J9").Value = 'M4 × 0.7 mm'
$ws.Range("K9").Value = '18-8 Stainless Steel'
$ws.Range("A10").Value = '16'
$ws.Range("B10").Value = 'Fully Threaded'
$ws.Range("C10").Value = '7.60'
$ws.Range("D10").Value = '2.20'
$ws.Range("E10").Value = '2.5'
$ws.Range("F10").Value = '80,000'
$ws.Range("G10").Value = '5'
$ws.Range("H10").Value = '95942A118'
$ws.Range("I10").Value = '11.50'
$ws.Range("J10").Value = 'M4 × 0.7 mm'
$ws.Range("K10").Value = '18-8 Stainless Steel'
$ws.Range("A11").Value = '10'
$ws.Range("B11").Value = 'Fully Threaded'
$ws.Range("C11").Value = '9.50'
$ws.Range("D11").Value = '2.75'
$ws.Range("E11").Value = '3'
$ws.Range("F11").Value = '80,000'
$ws.Range("G11").Value = '5'
$ws.Range("H11").Value = '95942A119'
$ws.Range("I11").Value = '9.64'
$ws.Range("J11").Value = 'M5 × 0.8 mm'
$ws.Range("K11").Value = '18-8 Stainless Steel'
$ws.Range("A12").Value = '12'
$ws.Range("B12").Value = 'Fully Threaded'
$ws.Range("C12").Value = '9.50'
$ws.Range("D12").Value = '2.75'
$ws.Range("E12").Value = '3'
$ws.Range("F12").Value = '80,000'
$ws.Range("G12").Value = '5'
$ws.Range("H12").Value = '95942A121'
$ws.Range("I12").Value = '10.34'
$ws.Range("J12").Value = 'M5 × 0.8 mm'
$ws.Range("K12").Value = '18-8 Stainless Steel'
$ws.Range("A13").Value = '16'
$ws.Range("B13").Value = 'Fully Threaded'
$ws.Range("C13").Value = '9.50'
$ws.Range("D13").Value = '2.75'
$ws.Range("E13").Value = '3'
$ws.Range("F13").Value = '80,000'
$ws.Range("G13").Value = '5'
$ws.Range("H13").Value = '95942A122'
$ws.Range("I13").Value = '11.15'
$ws.Range("J13").Value = 'M5 × 0.8 mm'
$ws.Range("K13").Value = '18-8 Stainless Steel'
$ws.Range("A14").Value = '20'
$ws.Range("B14").Value = 'Fully Threaded'
$ws.Range("C14").Value = '9.50'
$ws.Range("D14").Value = '2.75'
$ws.Range("E14").Value = '3'
$ws.Range("F14").Value = '80,000'
$ws.Range("G14").Value = '5'
$ws.Range("H14").Value = '95942A123'
$ws.Range("I14").Value = '12.43'
$ws.Range("J14").Value = 'M5 × 0.8 mm'
$ws.Range("K14").Value = '18-8 Stainless Steel'
$ws.Range("A15").Value = '10'
$ws.Range("B15").Value = 'Fully Threaded'
$ws.Range("C15").Value = '10.50'
$ws.Range("D15").Value = '3.30'
$ws.Range("E15").Value = '4'
$ws.Range("F15").Value = '80,000'
$ws.Range("G15").Value = '5'
$ws.Range("H15").Value = '95942A124'
$ws.Range("I15").Value = '10.66'
$ws.Range("J15").Value = 'M6 × 1 mm'
$ws.Range("K15").Value = '18-8 Stainless Steel'
$ws.Range("A16").Value = '12'
$ws.Range("B16").Value = 'Fully Threaded'
$ws.Range("C16").Value = '10.50'
$ws.Range("D16").Value = '3.30'
$ws.Range("E16").Value = '4'
$ws.Range("F16").Value = '80,000'
$ws.Range("G16").Value = '5'
$ws.Range("H16").Value = '95942A125'
$ws.Range("I16").Value = '11.61'
$ws.Range("J16").Value = 'M6 × 1 mm'
$ws.Range("K16").Value = '18-8 Stainless Steel'
$ws.Range("A17").Value = '16'
$ws.Range("B17").Value = 'Fully Threaded'
$ws.Range("C17").Value = '10.50'
$ws.Range("D17").Value = '3.30'
$ws.Range("E17").Value = '4'
$ws.Range("F17").Value = '80,000'
$ws.Range("G17").Value = '5'
$ws.Range("H17").Value = '95942A126'
$ws.Range("I17").Value = '12.65'
$ws.Range("J17").Value = 'M6 × 1 mm'
$ws.Range("K17").Value = '18-8 Stainless Steel'
$ws.Range("A18").Value = '20'
$ws.Range("B18").Value = 'Fully Threaded'
$ws.Range("C18").Value = '10.50'
$ws.Range("D18").Value = '3.30'
$ws.Range("E18").Value = '4'
$ws.Range("F18").Value = '80,000'
$ws.Range("G18").Value = '5'
$ws.Range("H18").Value = '95942A127'
$ws.Range("I18").Value = '14.22'
$ws.Range("J18").Value = 'M6 × 1 mm'
$ws.Range("K18").Value = '18-8 Stainless Steel'

# Restore the plain/default style on the rewritten cells (also resets NumberFormat)
$textRange.Style = "Normal"
